$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# fix#5432: PCTO Stage - aggiungere sede di svolgimento PCTO
#
# The "Patto formativo studente" paragraph referred to a single placeholder
# {A_SEDE} twice (once for the legal seat of the hosting organisation, once
# for the seat where the internship actually takes place). Split it into two
# distinct placeholders:
#   - first occurrence  ({A_SEDE})  -> {A_SEDE_LEGALE}
#   - second occurrence ({A_SEDE})  -> {A_SEDE_SVOLGIMENTO}
# ---------------------------------------------------------------------------

# Locate the first occurrence of the "{A_SEDE}" placeholder (the "sede legale").
$legale = $d.Content
$foundLegale = $legale.Find.Execute("{A_SEDE}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundLegale) {
    throw "Could not find the first {A_SEDE} placeholder (sede legale)."
}
$legale.Text = "{A_SEDE_LEGALE}"

# Locate the next occurrence, searching from right after the text we just
# inserted through to the end of the document (the "sede di svolgimento").
$svolgimento = $d.Range($legale.End, $d.Content.End)
$foundSvolgimento = $svolgimento.Find.Execute("{A_SEDE}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundSvolgimento) {
    throw "Could not find the second {A_SEDE} placeholder (sede di svolgimento)."
}
$svolgimento.Text = "{A_SEDE_SVOLGIMENTO}"

Write-Host "Placeholders updated: sede legale -> {A_SEDE_LEGALE}, sede svolgimento -> {A_SEDE_SVOLGIMENTO}"
